$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "J" (Obs) column for two existing rows: blank placeholder "-" -> "providenciar NF"
$ws.Range("J5").Value = "providenciar NF"
$ws.Range("J6").Value = "providenciar NF"

# Row 21's Obs column previously said "localizar NF" - now reads "providenciar NF"
$ws.Range("J21").Value = "providenciar NF"

# Append a new row (22) for the drone purchase.
# Set the values/formulas first (on the still-blank row), then copy row 21's
# formatting down onto row 22 - pasting formats only avoids Excel "helpfully"
# re-interpreting the freshly-typed date serial and minting a duplicate
# custom number format.
$ws.Range("A22").Value = "REGMEL"
$ws.Range("B22").Value = 45499
$ws.Range("C22").Formula = '="000638"'
$ws.Range("D22").Value = "CUSTO"
$ws.Range("E22").Value = "DRONE DJI AIR 2S FLY MORE COMBO"
$ws.Range("F22").Value = "und"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 9790
$ws.Range("I22").Formula = "=G22*H22"

$ws.Range("A21:I21").Copy()
$ws.Range("A22:I22").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 22 has no "Obs" entry, unlike row 21.
$ws.Range("J22").ClearContents()

$ws.Range("M20").Select()
